$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "date" column (F) values: add 2 days to each existing date serial,
# matching the diff (44914->44916, 44913->44915, 44912->44914, 44911->44913,
# 44910->44912, 44909->44911) for rows 2 through 7.
$ws.Range("F2").Value = 44916
$ws.Range("F3").Value = 44915
$ws.Range("F4").Value = 44914
$ws.Range("F5").Value = 44913
$ws.Range("F6").Value = 44912
$ws.Range("F7").Value = 44911
